# Apply translation / content updates to the smartfridge resources sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Ressource / Détail"
$ws.Range("C1").Value = "Quantité"
$ws.Range("D1").Value = "Coût estimé (€)"

# Row 2 - Hardware: Fridge prototype unit
$ws.Range("B2").Value = "Frigo standard modifié (base prototype)"
$ws.Range("E2").Value = "Support physique pour le prototype caméra+IA"

# Row 3 - Hardware: Camera module
$ws.Range("B3").Value = "Caméras internes HD grand angle résistantes au froid"
$ws.Range("E3").Value = "Doit supporter condensation, basse température"

# Row 4 - Hardware: Jetson Nano
$ws.Range("B4").Value = "Module IA embarqué (Jetson Nano / Raspberry Pi 5)"
$ws.Range("E4").Value = "Exécute l'IA en local (edge AI)"

# Row 5 - Software: TensorFlow / PyTorch
$ws.Range("B5").Value = "Stack IA vision (PyTorch + YOLO + TensorRT/ONNX)"
$ws.Range("E5").Value = "Open source → coût licence nul"

# Row 6 - Software: React Native + API tools
$ws.Range("B6").Value = "App mobile React Native + backend FastAPI/PostgreSQL"
$ws.Range("E6").Value = "Dev logiciel interne, pas de licence payante"

# Row 7 - Human: AI Engineer
$ws.Range("A7").Value = "Humain"
$ws.Range("B7").Value = "IA Engineer / Computer Vision"
$ws.Range("E7").Value = "1 mois d'IA/vision temps plein (estimation)"

# Row 8 - Human: IoT Engineer
$ws.Range("A8").Value = "Humain"
$ws.Range("B8").Value = "IoT / Embedded Engineer"
$ws.Range("E8").Value = "Intégration hardware + réseau Wi-Fi sécurisé"

# Row 9 - Human: Mobile Developer
$ws.Range("A9").Value = "Humain"
$ws.Range("B9").Value = "Mobile + Backend Developer"
$ws.Range("E9").Value = "DEV app + backend + intégration Drive"

# Row 10 - Infrastructure: Cloud server
$ws.Range("A10").Value = "Infra"
$ws.Range("B10").Value = "Serveur cloud (OVH / AWS) pour API + DB"
$ws.Range("E10").Value = "Hébergement / stockage 12 mois"

# Row 11 - Other: Training materials
$ws.Range("A11").Value = "Autre"
$ws.Range("B11").Value = "Temps réunions / tests utilisateurs pilotes"
$ws.Range("E11").Value = "Panels test familles pilotes / feedback terrain"
